$p = $ppt.ActivePresentation
$s = $p.Slides.Item(16)
$sh = $s.Shapes.Item(3)
$tbl = $sh.Table
Write-Host "Before:" $tbl.StyleId
$tbl.ApplyStyle("{13307E3E-D507-42F5-B440-E41BE3E34FB1}")
Write-Host "After:" $tbl.StyleId
